$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: new entry for "Valid Inagram" (LeetCode #242) ---
# Seed formatting on each cell by copying from a same-styled existing cell,
# then overwrite with the real value (keeps the shared style indices intact).

# A13 - No. (blue fill, like A5/A7/A8/A10/A12)
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 242

# B13 - Name (default style)
$ws.Range("B13").Value = "Valid Inagram"

# C13 - Main Topic (default style)
$ws.Range("C13").Value = "String"

# D13 - Tags (wrap text style, like D2/E2/I2)
$ws.Range("D2").Copy($ws.Range("D13"))
$ws.Range("D13").Value = "String, Hash table, Sorting"

# E13 - Related Topics (wrap text style)
$ws.Range("E2").Copy($ws.Range("E13"))
$ws.Range("E13").Value = "Dict"

# F13 - Level (default style)
$ws.Range("F13").Value = "Easy"

# G13 - Solutions count (default style)
$ws.Range("G13").Value = 2

# H13 - Status (red-font style, like H2..H12)
$ws.Range("H2").Copy($ws.Range("H13"))
$ws.Range("H13").Value = "✅"

# I13 - Understanding (wrap text style)
$ws.Range("I2").Copy($ws.Range("I13"))
$ws.Range("I13").Value = "Given 1 solution and copied 1 from solutions"

# Row height for the new entry
$ws.Rows.Item(13).RowHeight = 30

# --- Update the view state to match where the author left the cursor ---
[void]$ws.Range("I14").Select()
